# Cotações atualizadas - 2025-10-05
# Append the new daily quotation row (row 31) below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: date serial 45935 => 2025-10-05, same style/number format as
# the rows above it (yyyy-mm-dd hh:mm:ss).
$ws.Range("A31").Value = 45935
$ws.Range("A31").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Columns B-E: quotation values, stored as text using comma decimal
# separators, matching the existing data in the sheet.
$ws.Range("B31").Value = "21,4463"
$ws.Range("C31").Value = "15,2675"
$ws.Range("D31").Value = "15,4193"
$ws.Range("E31").Value = "15,4193"
